# Weekly update: insert a new price entry (Primera/Segunda) for
# "Zapallo italiano" - Agrícola del Norte S.A. de Arica, pushing the
# existing history down by one date (two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data (row 93/94), shifting
# all existing data rows (old 93-222) down to 95-224.
$ws.Rows("93:94").Insert()

# New "Primera" quality row for the latest date.
$ws.Range("A93").Value = 1
$ws.Range("B93").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C93").Value = "Arica y Parinacota"
$ws.Range("D93").Value = 44495
$ws.Range("E93").Value = 15
$ws.Range("F93").Value = 100112032
$ws.Range("G93").Value = "Zapallo italiano"
$ws.Range("H93").Value = "Huracán"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 120
$ws.Range("K93").Value = 7000
$ws.Range("L93").Value = 8000
$ws.Range("M93").Value = 7500
$ws.Range("N93").Value = "$/caja 70 unidades"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 107
$ws.Range("Q93").Value = 70
$ws.Range("R93").Value = "Hortaliza"

# New "Segunda" quality row for the latest date.
$ws.Range("A94").Value = 1
$ws.Range("B94").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C94").Value = "Arica y Parinacota"
$ws.Range("D94").Value = 44495
$ws.Range("E94").Value = 15
$ws.Range("F94").Value = 100112032
$ws.Range("G94").Value = "Zapallo italiano"
$ws.Range("H94").Value = "Huracán"
$ws.Range("I94").Value = "Segunda"
$ws.Range("J94").Value = 130
$ws.Range("K94").Value = 5000
$ws.Range("L94").Value = 6000
$ws.Range("M94").Value = 5500
$ws.Range("N94").Value = "$/caja 100 unidades"
$ws.Range("O94").Value = "Región de Arica y Parinacota"
$ws.Range("P94").Value = 55
$ws.Range("Q94").Value = 100
$ws.Range("R94").Value = "Hortaliza"
